# Update with restock suggestion
# - Forecast Comparison sheet: fill in Week_Start_Date (col B), refresh
#   Inventory Coverage / Stockout Risk / Reorder Urgency / Seasonality Index,
#   and collapse "Sales Volume Rank" + "Lifecycle Stage" into a single
#   "Lifecycle Stage" (now reporting "Decline") column.
# - Summary sheet: Max/Min Forecast Week become "N/A".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# --- Remove the "Sales Volume Rank" column (Q). This shifts the old
#     "Lifecycle Stage" column (R) left into Q, matching the new header
#     layout and the new dimension (A1:Q17 instead of A1:R17).
$ws.Range("Q1").EntireColumn.Delete()

# Make sure the Week_Start_Date column stays plain text (these values look
# like dates, but the source data keeps them as literal strings).
$ws.Range("B2:B17").NumberFormat = "@"

$rows = @(
    @{ Row = 2;  B = "2025-02-02"; L = 10.82;              M = "Low";  N = "Normal"; P = 1.12 },
    @{ Row = 3;  B = "2025-02-09"; L = 9.75;               M = "Low";  N = "Normal"; P = 1.17 },
    @{ Row = 4;  B = "2025-02-16"; L = 8.720000000000001;  M = "Low";  N = "Normal"; P = 1.12 },
    @{ Row = 5;  B = "2025-02-23"; L = 7.72;               M = "Low";  N = "Normal"; P = 1.1  },
    @{ Row = 6;  B = "2025-03-02"; L = 6.66;               M = "Low";  N = "Normal"; P = 0.84 },
    @{ Row = 7;  B = "2025-03-09"; L = 5.71;               M = "Low";  N = "Normal"; P = 0.9  },
    @{ Row = 8;  B = "2025-03-16"; L = 4.78;               M = "Low";  N = "Normal"; P = 1.1  },
    @{ Row = 9;  B = "2025-03-23"; L = 3.77;               M = "Low";  N = "Normal"; P = 1.08 },
    @{ Row = 10; B = "2025-03-30"; L = 2.77;               M = "Low";  N = "Normal"; P = 1.13 },
    @{ Row = 11; B = "2025-04-06"; L = 1.8;                M = "Low";  N = "Normal"; P = 1.04 },
    @{ Row = 12; B = "2025-04-13"; L = 0.8100000000000001; M = "Low";  N = "Urgent"; P = 1.18 },
    @{ Row = 13; B = "2025-04-20"; L = 0;                  M = "High"; N = "Urgent"; P = 0.93 },
    @{ Row = 14; B = "2025-04-27"; L = 0;                  M = "High"; N = "Urgent"; P = 1.15 },
    @{ Row = 15; B = "2025-05-04"; L = 0;                  M = "High"; N = "Urgent"; P = 0.92 },
    @{ Row = 16; B = "2025-05-11"; L = 0;                  M = "High"; N = "Urgent"; P = 0.97 },
    @{ Row = 17; B = "2025-05-18"; L = 0;                  M = "High"; N = "Urgent"; P = 1.17 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B    # Week_Start_Date
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # Inventory Coverage
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # Stockout Risk
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # Reorder Urgency
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # Seasonality Index
}

# New merged "Lifecycle Stage" column (previously R, now Q after the delete)
$ws.Range("Q2:Q17").Value = "Decline"

# --- Summary sheet: forecast week extremes are no longer available.
$summary.Range("B13").Value = "N/A"
$summary.Range("B15").Value = "N/A"
